$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AddOpportunity" (sheet1): append 3 rows (11-13) cloned from row 2,
# with a new JobType value (col C) for each, styled with wrap+vcenter.
# ---------------------------------------------------------------------------
$wsOpp = $wb.Worksheets.Item("AddOpportunity")

$oppTemplate = @{
    1  = "Techno Coatings, Inc."
    2  = "Techno Coatings, Inc."
    4  = "BUS - Business Services"
    5  = "Dealership & Rental Services"
    6  = "No"
    7  = "No"
    8  = "Accountant"
    9  = "No"
    10 = "No"
    11 = "AM"
    12 = "HL Capital, Inc."
    13 = "Do Not Disclose"
    14 = "Emre Abale"
    15 = "10"
    16 = "10"
    17 = "10"
    18 = "Public Equity"
    19 = "Public Equity"
    20 = "9999"
    21 = "Test"
    22 = "Chris Lord"
    23 = "Yes, separate signed agreement"
    24 = "Cleared"
    25 = "CF"
    26 = "Consulting"
    27 = "10"
    28 = "10"
    29 = "Emre Abale"
    30 = "Yes"
}

$oppNewJobTypes = @{
    11 = "Strategy"
    12 = "Post Merger Integration"
    13 = "Valuation Advisory"
}

foreach ($row in 11, 12, 13) {
    foreach ($col in $oppTemplate.Keys) {
        $wsOpp.Cells.Item($row, $col).Value = $oppTemplate[$col]
    }
    $wsOpp.Cells.Item($row, 3).Value = $oppNewJobTypes[$row]
}

# Apply the wrap-text / vertically-centered style to the new JobType cells.
# Build the style directly on C11 (single cell -> no orphan intermediate xf),
# then copy the format (only) onto C12:C13 so every cell lands on the very
# same cellXfs entry.
$wsOpp.Cells.Item(11, 3).VerticalAlignment = -4108
$wsOpp.Cells.Item(11, 3).WrapText = $true
$wsOpp.Range("C11").Copy()
$wsOpp.Range("C12:C13").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet "Engagement" (sheet4): append 3 rows (11-13) cloned from row 4,
# with a new UpdateJobType value (col D) for each, all four cells styled
# with wrap+vcenter.
# ---------------------------------------------------------------------------
$wsEng = $wb.Worksheets.Item("Engagement")

$engNewJobTypes = @{
    11 = "Strategy"
    12 = "Post Merger Integration"
    13 = "Valuation Advisory"
}

foreach ($row in 11, 12, 13) {
    $wsEng.Cells.Item($row, 1).Value = "Advisory"
    $wsEng.Cells.Item($row, 2).Value = "Advisory (CF)"
    $wsEng.Cells.Item($row, 3).Value = "HL Capital, Inc."
    $wsEng.Cells.Item($row, 4).Value = $engNewJobTypes[$row]
}

$wsEng.Cells.Item(11, 1).VerticalAlignment = -4108
$wsEng.Cells.Item(11, 1).WrapText = $true
$wsEng.Range("A11").Copy()
$wsEng.Range("B11:D11").PasteSpecial(-4122)
$wsEng.Range("A12:D13").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# View state: AddOpportunity keeps a selection but is no longer the active
# tab; Engagement becomes the active tab with its own selection.
# ---------------------------------------------------------------------------
$wsOpp.Range("C11:C13").Select()
$wsEng.Activate()
$wsEng.Range("D16").Select()
